# REVER_DailyTracker update (per commit "Add files via upload"):
# Fill in the "Ram" worksheet's daily log for 2020-10-06 .. 2020-10-15
# (rows 8-17), then leave the "Ram" tab active/selected the way the
# author left it (D20:D21 selected) instead of "Monisha".

$wb = $excel.ActiveWorkbook
$ram = $wb.Worksheets.Item("Ram")
$donorTask = $wb.Worksheets.Item("Sabeena")   # has the A/B/C/D/E "task row" style combo already

# ---- helper donor ranges already present in the workbook ----
# Sabeena!A6:E6  -> styles: No=1, Date=31, Application=29, Task=30, %=7
# Ram!A5:E5      -> styles: No=1, Date=31, Application=29, Task(Week off)=36, %=41 (Week-off template)
# Ram!F5         -> style 42 (Week-off "F" blank cell)
# Ram!B23        -> style 43 (WIP fill, centered) - reuse for F column WIP cells

function Copy-Format($srcRange, $dstRange) {
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats
}

# --- Row 8: No=6, 2020-10-06, Soniya, creating Setup file, 100%, Completed ---
Copy-Format $donorTask.Range("A6:E6") $ram.Range("A8:E8")
$ram.Cells.Item(8,1).Value = 6
$ram.Cells.Item(8,2).Value = 44110
$ram.Cells.Item(8,3).Value = "Soniya"
$ram.Cells.Item(8,4).Value = "creating Setup file"
$ram.Cells.Item(8,5).Value = 1

$f8 = $ram.Cells.Item(8,6)
$f8.Value = "Completed"
$f8.Borders.ColorIndex = 1
$f8.Borders.LineStyle = 1
$f8.Borders.Weight = 2
$f8.Interior.Color = 5287936
$f8.HorizontalAlignment = -4108

# --- Row 9: No=7, 2020-10-07, QMVAR, Layout Designing, (blank %), WIP ---
Copy-Format $donorTask.Range("A6:E6") $ram.Range("A9:E9")
$ram.Cells.Item(9,1).Value = 7
$ram.Cells.Item(9,2).Value = 44111
$ram.Cells.Item(9,3).Value = "QMVAR"
$ram.Cells.Item(9,4).Value = "Layout Designing"

Copy-Format $ram.Range("B23") $ram.Range("F9")
$ram.Cells.Item(9,6).Value = "WIP"

# --- Row 10: No=8, 2020-10-08, Emplogin, Adding assests Detatils, 50%, WIP ---
Copy-Format $donorTask.Range("A6:E6") $ram.Range("A10:E10")
$ram.Cells.Item(10,1).Value = 8
$ram.Cells.Item(10,2).Value = 44112
$ram.Cells.Item(10,3).Value = "Emplogin"
$ram.Cells.Item(10,4).Value = "Adding assests Detatils"
$ram.Cells.Item(10,5).Value = 0.5

Copy-Format $ram.Range("B23") $ram.Range("F10")
$ram.Cells.Item(10,6).Value = "WIP"

# --- Row 11: No=9, 2020-10-09, Emplogin, Adding assests Detatils, 100%, Completed ---
Copy-Format $donorTask.Range("A6:E6") $ram.Range("A11:E11")
$ram.Cells.Item(11,1).Value = 9
$ram.Cells.Item(11,2).Value = 44113
$ram.Cells.Item(11,3).Value = "Emplogin"
$ram.Cells.Item(11,4).Value = "Adding assests Detatils"
$ram.Cells.Item(11,5).Value = 1

Copy-Format $ram.Range("F8") $ram.Range("F11")
$ram.Cells.Item(11,6).Value = "Completed"

# --- Row 12: No=10, 2020-10-10, Week off ---
Copy-Format $ram.Range("A5:E5") $ram.Range("A12:E12")
$ram.Cells.Item(12,1).Value = 10
$ram.Cells.Item(12,2).Value = 44114
$ram.Cells.Item(12,4).Value = "Week off"

Copy-Format $ram.Range("F5") $ram.Range("F12")

# --- Row 13: No=11, 2020-10-11, Week off ---
Copy-Format $ram.Range("A5:E5") $ram.Range("A13:E13")
$ram.Cells.Item(13,1).Value = 11
$ram.Cells.Item(13,2).Value = 44115
$ram.Cells.Item(13,4).Value = "Week off"

Copy-Format $ram.Range("F5") $ram.Range("F13")

# --- Row 14: No=12, 2020-10-12, Emplogin, Adding assests Detatils and responisble view, 100%, Completed ---
Copy-Format $donorTask.Range("A6:E6") $ram.Range("A14:E14")
$ram.Cells.Item(14,1).Value = 12
$ram.Cells.Item(14,2).Value = 44116
$ram.Cells.Item(14,3).Value = "Emplogin"
$ram.Cells.Item(14,4).Value = "Adding assests Detatils and responisble view"
$ram.Cells.Item(14,5).Value = 1

Copy-Format $ram.Range("F8") $ram.Range("F14")
$ram.Cells.Item(14,6).Value = "Completed"

# --- Row 15: No=13, 2020-10-13, QMVAR, Adding assests Detatils, 100%, WIP ---
Copy-Format $donorTask.Range("A6:E6") $ram.Range("A15:E15")
$ram.Cells.Item(15,1).Value = 13
$ram.Cells.Item(15,2).Value = 44117
$ram.Cells.Item(15,3).Value = "QMVAR"
$ram.Cells.Item(15,4).Value = "Adding assests Detatils"
$ram.Cells.Item(15,5).Value = 1

Copy-Format $ram.Range("B23") $ram.Range("F15")
$ram.Cells.Item(15,6).Value = "WIP"

# --- Row 16: No=14, 2020-10-14, QMVAR, stored management added in Monthly target, 100%, Completed ---
Copy-Format $donorTask.Range("A6:E6") $ram.Range("A16:E16")
$ram.Cells.Item(16,1).Value = 14
$ram.Cells.Item(16,2).Value = 44118
$ram.Cells.Item(16,3).Value = "QMVAR"
$ram.Cells.Item(16,4).Value = "stored management added in Monthly target"
$ram.Cells.Item(16,5).Value = 1

Copy-Format $ram.Range("F8") $ram.Range("F16")
$ram.Cells.Item(16,6).Value = "Completed"

# --- Row 17: No=15, 2020-10-15, QMVAR, adding layout in analysis file upload, 100%, WIP ---
Copy-Format $donorTask.Range("A6:E6") $ram.Range("A17:E17")
$ram.Cells.Item(17,1).Value = 15
$ram.Cells.Item(17,2).Value = 44119
$ram.Cells.Item(17,3).Value = "QMVAR"
$ram.Cells.Item(17,4).Value = "adding layout in analysis file upload"
$ram.Cells.Item(17,5).Value = 1

Copy-Format $ram.Range("B23") $ram.Range("F17")
$ram.Cells.Item(17,6).Value = "WIP"

# --- Activate Ram tab and restore the author's last selection there ---
$ram.Activate()
$ram.Range("D20:D21").Select()
